$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 48.0404218882511
$ws.Range("B1").Value = 97.5989306800062
$ws.Range("C1").Value = 184.821959577884
$ws.Range("D1").Value = 161.672081780467
$ws.Range("E1").Value = 26.0244057634959
$ws.Range("F1").Value = 9.45886429839715
$ws.Range("G1").Value = 198.601725324337
$ws.Range("H1").Value = 143.573255996952
$ws.Range("I1").Value = 38.2128565750145
$ws.Range("J1").Value = 162.330589053375
$ws.Range("A2").Value = 104.984851137262
$ws.Range("B2").Value = 23.8904076739635
$ws.Range("C2").Value = 60.0762073230353
$ws.Range("D2").Value = 77.2711618231941
$ws.Range("E2").Value = 149.651743355045
$ws.Range("F2").Value = 92.0145681556382
$ws.Range("G2").Value = 15.5266159286381
$ws.Range("H2").Value = 19.3194829948803
$ws.Range("I2").Value = 6.53834278068428
$ws.Range("J2").Value = 89.9085517460054
$ws.Range("A3").Value = 37.1586835184873
$ws.Range("B3").Value = 27.5334047281805
$ws.Range("C3").Value = 151.242737356221
$ws.Range("D3").Value = 134.950748987007
$ws.Range("E3").Value = 159.091253559613
$ws.Range("F3").Value = 45.4644146587068
$ws.Range("G3").Value = 72.0924906768335
$ws.Range("H3").Value = 194.582911205749
$ws.Range("I3").Value = 122.188384329057
$ws.Range("J3").Value = 199.384425487083
$ws.Range("A4").Value = 33.3120031437427
$ws.Range("B4").Value = 105.305086497825
$ws.Range("C4").Value = 76.8719453722574
$ws.Range("D4").Value = 1.45210940458444
$ws.Range("E4").Value = 176.642692450733
$ws.Range("F4").Value = 121.071253400795
$ws.Range("G4").Value = 175.423231150686
$ws.Range("H4").Value = 31.9730667546266
$ws.Range("I4").Value = 111.062208242278
$ws.Range("J4").Value = 14.4631481796797
$ws.Range("A5").Value = 120.360910203476
$ws.Range("B5").Value = 76.5528200550716
$ws.Range("C5").Value = 50.061680213577
$ws.Range("D5").Value = 15.6393317578544
$ws.Range("E5").Value = 28.2757044901493
$ws.Range("F5").Value = 134.180644310164
$ws.Range("G5").Value = 73.6901930876496
$ws.Range("H5").Value = 9.35079846966583
$ws.Range("I5").Value = 172.538007037965
$ws.Range("J5").Value = 147.579582942454
$ws.Range("A6").Value = 104.525548361487
$ws.Range("B6").Value = 88.2553300299939
$ws.Range("C6").Value = 25.9046786585379
$ws.Range("D6").Value = 178.122212355082
$ws.Range("E6").Value = 45.008178448774
$ws.Range("F6").Value = 20.5070171600706
$ws.Range("G6").Value = 146.356193323786
$ws.Range("H6").Value = 49.8712105908763
$ws.Range("I6").Value = 2.58082822085397
$ws.Range("J6").Value = 180.559991104789
$ws.Range("A7").Value = 137.366373621564
$ws.Range("B7").Value = 4.01881411858779
$ws.Range("C7").Value = 21.3848716678959
$ws.Range("D7").Value = 38.8284310879318
$ws.Range("E7").Value = 129.018585909632
$ws.Range("F7").Value = 199.679764639437
$ws.Range("G7").Value = 147.018462301706
$ws.Range("H7").Value = 58.6240979184509
$ws.Range("I7").Value = 100.628469372461
$ws.Range("J7").Value = 28.5804899542502
$ws.Range("A8").Value = 116.591337004952
$ws.Range("B8").Value = 183.553549174011
$ws.Range("C8").Value = 108.257274752603
$ws.Range("D8").Value = 192.075194601005
$ws.Range("E8").Value = 169.54764154253
$ws.Range("F8").Value = 160.605863463416
$ws.Range("G8").Value = 177.471724514603
$ws.Range("H8").Value = 135.603405598366
$ws.Range("I8").Value = 106.675044496858
$ws.Range("J8").Value = 24.9106092494496
$ws.Range("A9").Value = 171.774221571057
$ws.Range("B9").Value = 62.7416922071677
$ws.Range("C9").Value = 22.0449041677848
$ws.Range("D9").Value = 174.608801386603
$ws.Range("E9").Value = 94.8588771255961
$ws.Range("F9").Value = 145.056673113749
$ws.Range("G9").Value = 79.4004078392873
$ws.Range("H9").Value = 98.7497330171753
$ws.Range("I9").Value = 156.44393095581
$ws.Range("J9").Value = 156.135675290663
$ws.Range("A10").Value = 174.715060077009
$ws.Range("B10").Value = 125.55202055981
$ws.Range("C10").Value = 29.3922385337726
$ws.Range("D10").Value = 130.502217137489
$ws.Range("E10").Value = 77.0967745581161
$ws.Range("F10").Value = 116.342096084888
$ws.Range("G10").Value = 55.1679483871758
$ws.Range("H10").Value = 11.2332491256451
$ws.Range("I10").Value = 86.620745848222
$ws.Range("J10").Value = 28.5959398507122
$ws.Range("A11").Value = 187.162182008457
$ws.Range("B11").Value = 15.0660951691987
$ws.Range("C11").Value = 108.722329097205
$ws.Range("D11").Value = 143.957517083714
$ws.Range("E11").Value = 30.9882459375021
$ws.Range("F11").Value = 120.971999187475
$ws.Range("G11").Value = 179.998055277391
$ws.Range("H11").Value = 33.8294840575333
$ws.Range("I11").Value = 8.5745708125525
$ws.Range("J11").Value = 84.4023149853583
$ws.Range("A12").Value = 43.0352926454671
$ws.Range("B12").Value = 10.7527877254192
$ws.Range("C12").Value = 143.196166094018
$ws.Range("D12").Value = 177.670218971404
$ws.Range("E12").Value = 8.78576953373187
$ws.Range("F12").Value = 74.6246814143959
$ws.Range("G12").Value = 181.973909950803
$ws.Range("H12").Value = 46.7760702812933
$ws.Range("I12").Value = 143.969553962336
$ws.Range("J12").Value = 183.961912795884
$ws.Range("A13").Value = 120.27935680015
$ws.Range("B13").Value = 48.2687292845308
$ws.Range("C13").Value = 102.18016696264
$ws.Range("D13").Value = 144.492794081798
$ws.Range("E13").Value = 53.8654298772409
$ws.Range("F13").Value = 191.039316445142
$ws.Range("G13").Value = 154.161310640239
$ws.Range("H13").Value = 177.755057615114
$ws.Range("I13").Value = 114.978420042888
$ws.Range("J13").Value = 53.2055454576414
$ws.Range("A14").Value = 105.43791507624
$ws.Range("B14").Value = 166.238475388958
$ws.Range("C14").Value = 48.9826597501443
$ws.Range("D14").Value = 78.0791046461459
$ws.Range("E14").Value = 37.7484272409922
$ws.Range("F14").Value = 156.708126401858
$ws.Range("G14").Value = 154.019363109963
$ws.Range("H14").Value = 78.0873870840703
$ws.Range("I14").Value = 143.6205554491
$ws.Range("J14").Value = 173.886877938121
$ws.Range("A15").Value = 165.058617836357
$ws.Range("B15").Value = 45.5709237817539
$ws.Range("C15").Value = 90.1751622046228
$ws.Range("D15").Value = 72.0416159704521
$ws.Range("E15").Value = 113.100382645196
$ws.Range("F15").Value = 163.96227235159
$ws.Range("G15").Value = 182.355854465792
$ws.Range("H15").Value = 51.7220195623683
$ws.Range("I15").Value = 121.716447603757
$ws.Range("J15").Value = 2.47209314372022
$ws.Range("A16").Value = 134.368186134085
$ws.Range("B16").Value = 8.3918781058825
$ws.Range("C16").Value = 67.2636951633094
$ws.Range("D16").Value = 102.658833052338
$ws.Range("E16").Value = 108.316583050562
$ws.Range("F16").Value = 138.893452723927
$ws.Range("G16").Value = 112.885928206558
$ws.Range("H16").Value = 164.229535015407
$ws.Range("I16").Value = 90.0920872064736
$ws.Range("J16").Value = 39.9489294923604
$ws.Range("A17").Value = 166.810688547236
$ws.Range("B17").Value = 2.24299766227743
$ws.Range("C17").Value = 118.851064014645
$ws.Range("D17").Value = 155.369025354911
$ws.Range("E17").Value = 178.964399909118
$ws.Range("F17").Value = 76.7968172565088
$ws.Range("G17").Value = 161.770127975275
$ws.Range("H17").Value = 65.1170614478723
$ws.Range("I17").Value = 139.921791730412
$ws.Range("J17").Value = 52.0776431318734
$ws.Range("A18").Value = 120.605318304433
$ws.Range("B18").Value = 103.886522866733
$ws.Range("C18").Value = 103.155514832193
$ws.Range("D18").Value = 170.082676024215
$ws.Range("E18").Value = 18.9032949595262
$ws.Range("F18").Value = 74.7084330183959
$ws.Range("G18").Value = 158.093567079908
$ws.Range("H18").Value = 30.1385509921883
$ws.Range("I18").Value = 31.3924114366027
$ws.Range("J18").Value = 89.9031575256508
$ws.Range("A19").Value = 8.68346197934982
$ws.Range("B19").Value = 102.439291077871
$ws.Range("C19").Value = 56.0386100113572
$ws.Range("D19").Value = 112.506326899168
$ws.Range("E19").Value = 118.837359323556
$ws.Range("F19").Value = 97.0460369703574
$ws.Range("G19").Value = 98.9747802256489
$ws.Range("H19").Value = 146.323826697806
$ws.Range("I19").Value = 169.762521595584
$ws.Range("J19").Value = 98.8549745170656
$ws.Range("A20").Value = 43.8221981953002
$ws.Range("B20").Value = 189.789828094556
$ws.Range("C20").Value = 187.995299877597
$ws.Range("D20").Value = 103.67162595674
$ws.Range("E20").Value = 7.0761893908848
$ws.Range("F20").Value = 162.81562017408
$ws.Range("G20").Value = 126.719859767109
$ws.Range("H20").Value = 134.806136849712
$ws.Range("I20").Value = 93.0772160613338
$ws.Range("J20").Value = 36.3035653886868
